# edit.ps1 - Apply "Add PA work and some edits to CHN position" changes
#
# Targets the four Planned Parenthood Clinical Health Network (CHN) bullet
# paragraphs (the *first* occurrence of this job block in the resume - the
# Betterhealth job further down the document has very similar-looking text
# and must not be touched), plus a small cosmetic cleanup on the
# "Betterhealth" company-name run.

$d = $word.ActiveDocument

function Get-ParaIndex($containsText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $containsText + "*")) {
            return $i
        }
    }
    return -1
}

# --- 1. "Configure, manage, validate, test, troubleshoot..." bullet ---
$oldText1 = "Configure, manage, validate, test, troubleshoot, and maintain AWS Linux and Windows servers, storage, database, and network resources using direct configuration and abstraction tools, including config management, automation, and orchestration with Terraform, CloudFormation."
$idx1 = Get-ParaIndex($oldText1)
$r1 = $d.Paragraphs.Item($idx1).Range
$r1.Find.Execute($oldText1, $true, $true, $false, $false, $false, $true, 1, $false, `
    "Configure, manage, validate, test, troubleshoot, and maintain AWS Linux and Windows servers, storage, database, and network resources using direct configuration and abstraction automation, and orchestration tools such as Terraform and CloudFormation.", `
    2) | Out-Null

# --- 2. "Provide 3rd and 4th tier support..." bullet becomes the new      ---
#        "security and compliance tooling" bullet (PA / security work)    ---
$oldText2 = "Provide 3rd and 4th tier support to Service Desk and Planned Parenthood affiliates to resolve complex networking and device problems in a multi-region, multi-gigabit-per-second AWS and on premises network environment."
$idx2 = Get-ParaIndex($oldText2)
$r2 = $d.Paragraphs.Item($idx2).Range
$r2.Find.Execute($oldText2, $true, $true, $false, $false, $false, $true, 1, $false, `
    "Configure, manage, validate, test, and troubleshoot implementation of new security and compliance tooling in partnership with Information Security Team. Recent experience with Palo Alto Cortex XDR, CASB Cloud Security, and SaaS Security Posture Management. ", `
    2) | Out-Null

# --- 3. "Onboard new vendors..." bullet becomes the (reworded / expanded) ---
#        "Provide 3rd and 4th tier support..." bullet                     ---
$oldText3 = "Onboard new vendors and affiliates into the environment, including setting up site-to-site VPN tunnels through Palo Alto Prisma SD-WAN, creating Okta SAML sign on integrations, user and group management, and peripheral device configuration assistance."
$idx3 = Get-ParaIndex($oldText3)
$r3 = $d.Paragraphs.Item($idx3).Range
$r3.Find.Execute($oldText3, $true, $true, $false, $false, $false, $true, 1, $false, `
    "Provide 3rd and 4th tier support to Service Desk and Planned Parenthood affiliates to resolve complex networking and device problems in a multi-region, multi-gigabit-per-second hybrid AWS and on-premise network environment. Issues include site-to-site IPSec VPN tunnels setup or troubleshooting, SAML sign on integrations, wireless network issues, and peripheral device configuration.", `
    2) | Out-Null

# --- 4. "Work cross-functionally..." bullet gains "security and          ---
#        compliance," in the middle                                      ---
$oldText4 = "Work cross-functionally with other teams to develop and maintain solutions to support patient care, business services, and ensure appropriate use of IT resources across the organization."
$idx4 = Get-ParaIndex($oldText4)
$r4 = $d.Paragraphs.Item($idx4).Range
$r4.Find.Execute($oldText4, $true, $true, $false, $false, $false, $true, 1, $false, `
    "Work cross-functionally with other teams to develop and maintain solutions to support patient care, business services, security and compliance, and ensure appropriate use of IT resources across the organization.", `
    2) | Out-Null

# --- 5. "Betterhealth" company line: merge the "Betterhealth" run with   ---
#        the " A Planned Parenthood Partnership, Remote," run into a      ---
#        single run (drops the spell-check split the two used to sit in) ---
$idx5 = Get-ParaIndex("Betterhealth A Planned Parenthood Partnership, Remote,")
$r5 = $d.Paragraphs.Item($idx5).Range
$r5.Find.Execute("Betterhealth A Planned Parenthood Partnership, Remote,", $true, $true, $false, $false, $false, $true, 1, $false, `
    "Betterhealth A Planned Parenthood Partnership, Remote,", `
    2) | Out-Null

Write-Host "Paragraph indices touched:" $idx1 $idx2 $idx3 $idx4 $idx5
